# Update the cryptos price list (columns D = Price, E = Volume(1h)) with the
# latest scraped values. Cells that look like plain decimal numbers (e.g.
# "591.44") are written with a leading apostrophe so Excel keeps them stored
# as text, matching the original workbook's text-based Price column.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.073.10"
$ws.Range("E2").Value = "  -0.99%  "
$ws.Range("D3").Value = "2.603.02"
$ws.Range("E3").Value = "  -0.23%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'591.44"
$ws.Range("E5").Value = "  -0.85%  "
$ws.Range("D6").Value = "'151.67"
$ws.Range("E6").Value = "  -2.13%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("E8").Value = "  +0.43%  "
$ws.Range("D9").Value = "2.599.81"
$ws.Range("E9").Value = "  -0.27%  "
$ws.Range("E10").Value = "  -3.48%  "
$ws.Range("E11").Value = "  +0.25%  "
$ws.Range("E12").Value = "  -1.68%  "
$ws.Range("E13").Value = "  -2.76%  "
$ws.Range("E14").Value = "  -0.59%  "
$ws.Range("D15").Value = "3.074.55"
$ws.Range("E15").Value = "  -0.17%  "
$ws.Range("E16").Value = "  -4.40%  "
$ws.Range("D17").Value = "66.825.10"
$ws.Range("E17").Value = "  -1.10%  "
$ws.Range("D18").Value = "2.600.51"
$ws.Range("E18").Value = "  -0.45%  "
$ws.Range("D19").Value = "'363.27"
$ws.Range("E19").Value = "  -0.62%  "
$ws.Range("E20").Value = "  -1.56%  "
$ws.Range("D21").Value = "'7.36"
$ws.Range("E21").Value = "  -4.87%  "
$ws.Range("E22").Value = "  -0.48%  "
$ws.Range("D23").Value = "'2.05"
$ws.Range("E23").Value = "  +0.94%  "
$ws.Range("E24").Value = "  -0.12%  "
$ws.Range("D25").Value = "'10.04"
$ws.Range("E25").Value = "  +1.35%  "
$ws.Range("D26").Value = "'67.40"
$ws.Range("E26").Value = "  -0.52%  "
$ws.Range("D27").Value = "2.731.45"
$ws.Range("E27").Value = "  -0.68%  "
$ws.Range("D28").Value = "'584.11"
$ws.Range("E28").Value = "  +0.64%  "
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  +0.41%  "
$ws.Range("E30").Value = "  -3.75%  "
$ws.Range("E31").Value = "  -3.52%  "
$ws.Range("E32").Value = "  -3.76%  "
$ws.Range("D33").Value = "'1.79"
$ws.Range("E33").Value = "  -3.50%  "
$ws.Range("E34").Value = "  -0.02%  "
$ws.Range("E35").Value = "  -9.04%  "
$ws.Range("E36").Value = "  -2.16%  "
$ws.Range("E37").Value = "  -2.20%  "
$ws.Range("D38").Value = "'156.24"
$ws.Range("E38").Value = "  +0.79%  "
$ws.Range("D39").Value = "'18.84"
$ws.Range("E39").Value = "  -2.96%  "
$ws.Range("D40").Value = "'0.366"
$ws.Range("E40").Value = "  -1.22%  "
$ws.Range("D41").Value = "'5.23"
$ws.Range("E41").Value = "  -2.84%  "
$ws.Range("D42").Value = "'1.80"
$ws.Range("E42").Value = "  -3.89%  "
$ws.Range("E43").Value = "  +0.55%  "
$ws.Range("E44").Value = "  +2.00%  "
$ws.Range("E45").Value = "  -0.03%  "
$ws.Range("D46").Value = "0.0₆0296"
$ws.Range("E46").Value = "  +0.96%  "
$ws.Range("D47").Value = "'153.87"
$ws.Range("E47").Value = "  -1.51%  "
$ws.Range("E48").Value = "  -0.65%  "
$ws.Range("E49").Value = "  +3.31%  "
$ws.Range("E50").Value = "  -2.38%  "
$ws.Range("E51").Value = "  -2.86%  "
